# Apply updated imputed statistics to the "Bayesian imputation" rows (16-22)
# in the covid19deaths_impute_age_group2 results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (age_group 18-29, Bayesian imputation)
$ws.Range("B16").Value = 3950.205
$ws.Range("D16").Value = 3810.975
$ws.Range("E16").Value = 4085
$ws.Range("G16").Value = 167.63

# Row 17 (age_group 30-39, Bayesian imputation)
$ws.Range("B17").Value = 8758.209
$ws.Range("D17").Value = 8566.975
$ws.Range("E17").Value = 8950.025

# Row 18 (age_group 40-49, Bayesian imputation)
$ws.Range("B18").Value = 18304.913
$ws.Range("D18").Value = 18068.975
$ws.Range("E18").Value = 18562.075
$ws.Range("G18").Value = 62.12

# Row 19 (age_group 50-64, Bayesian imputation)
$ws.Range("B19").Value = 65777.407
$ws.Range("D19").Value = 65426.975
$ws.Range("E19").Value = 66119.075

# Row 20 (age_group 65-74, Bayesian imputation)
$ws.Range("B20").Value = 90358.363
$ws.Range("D20").Value = 90011.9
$ws.Range("E20").Value = 90696.175
$ws.Range("G20").Value = 10.11

# Row 21 (age_group 75-84, Bayesian imputation)
$ws.Range("B21").Value = 113356.84
$ws.Range("D21").Value = 112999.9
$ws.Range("E21").Value = 113699
$ws.Range("G21").Value = 6.98

# Row 22 (age_group 85+, Bayesian imputation)
$ws.Range("B22").Value = 129195.39
$ws.Range("D22").Value = 128871.975
$ws.Range("E22").Value = 129540
$ws.Range("G22").Value = 5.48
